# Fruta / hortaliza, semanal
# Insert a new weekly record for "Vega Modelo de Temuco - Membrillo" at row 270,
# pushing the existing rows 270-304 down to 271-305.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 270 (shifts rows 270:304 -> 271:305,
# inheriting formatting/number-format from the row above, same as Excel's
# default Insert behaviour).
$ws.Rows.Item(270).Insert()

# Populate the newly inserted row 270 with this week's data.
$ws.Range("A270").Value = 10
$ws.Range("B270").Value = "Vega Modelo de Temuco"
$ws.Range("C270").Value = "La Araucanía"
$ws.Range("D270").Value = 45124
$ws.Range("E270").Value = 9
$ws.Range("F270").Value = "Fruta"
$ws.Range("G270").Value = 100104
$ws.Range("H270").Value = "Frutos de pepita"
$ws.Range("I270").Value = 100104003
$ws.Range("J270").Value = "Membrillo"
$ws.Range("K270").Value = "Champion"
$ws.Range("L270").Value = "Primera"
$ws.Range("M270").Value = 125
$ws.Range("N270").Value = 14000
$ws.Range("O270").Value = 14000
$ws.Range("P270").Value = 14000
$ws.Range("Q270").Value = "$/bandeja 18 kilos granel"
$ws.Range("R270").Value = "Región de O'Higgins"
$ws.Range("S270").Value = 778
$ws.Range("T270").Value = 18
